$wb = $excel.ActiveWorkbook

# --- Sheet1: restore its literal text (the source file's shared-strings
# part failed to round-trip, so re-assert the existing "A" label) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "A"

# --- add Sheet2 after Sheet1 and make it the active sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- A1: multi-run rich text "Test Weight (Kgs)" ---
$ws2.Range("A1").Value = "Test Weight (Kgs)"

$run1 = $ws2.Range("A1").Characters(1, 13)
$run1.Font.Bold = $true
$run1.Font.Size = 9
$run1.Font.Name = "Tahoma"
$run1.Font.Color = 16777215

$run2 = $ws2.Range("A1").Characters(14, 3)
$run2.Font.Bold = $true
$run2.Font.Size = 9
$run2.Font.Name = "Tahoma"
$run2.Font.Color = 16777215

$run3 = $ws2.Range("A1").Characters(17, 1)
$run3.Font.Bold = $true
$run3.Font.Size = 9
$run3.Font.Name = "Tahoma"
$run3.Font.Color = 16777215

# --- cell style: bold white Tahoma on blue fill, thin gray border, wrap ---
$cell = $ws2.Range("A1")
$cell.Font.Bold = $true
$cell.Font.Size = 9
$cell.Font.Name = "Tahoma"
$cell.Font.Color = 16777215
$cell.Interior.Color = 10177632
$cell.Interior.Pattern = 1
$cell.Borders.LineStyle = 1
$cell.Borders.Weight = 2
$cell.Borders.Color = 6908265
$cell.WrapText = $true
$cell.VerticalAlignment = -4160
$cell.RowHeight = 26

$ws2.Activate()

Write-Host "done"
